# Update the "想去人数" (want-to-go count) column F on the "展览" and
# "全部类型" worksheets with freshly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 2-7 updated) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 374
$wsExpo.Range("F3").Value = 107
$wsExpo.Range("F4").Value = 1579
$wsExpo.Range("F5").Value = 11
$wsExpo.Range("F6").Value = 23
$wsExpo.Range("F7").Value = 404

# --- Sheet "全部类型" (rows 4-9 updated) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1579
$wsAll.Range("F5").Value = 11
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 404
$wsAll.Range("F8").Value = 138
$wsAll.Range("F9").Value = 62
